$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 6-17: update the User-Story id / role / story / need text values.
# (cell styles for these rows are unchanged except D17, handled below)
# ---------------------------------------------------------------------------

$ws.Range("A6").Value  = "US001"
$ws.Range("D6").Value  = "mantenerse un stock de productos actualizado"

$ws.Range("A7").Value  = "US002"
$ws.Range("B7").Value  = "Usuario "
$ws.Range("C7").Value  = "Modificar una o varias características de un producto ya registrado en el sistema. "
$ws.Range("D7").Value  = "mantenerse un stock de productos actualizado"

$ws.Range("A8").Value  = "US003"
$ws.Range("C8").Value  = "Eliminar un producto existente."
$ws.Range("D8").Value  = "mantenerse un stock de productos actualizado"

$ws.Range("A9").Value  = "US004"
$ws.Range("C9").Value  = "Consultar un producto existente."
$ws.Range("D9").Value  = "buscar informacion relevante de un producto"

$ws.Range("A10").Value = "US005"
$ws.Range("C10").Value = "Registrar una venta"
$ws.Range("D10").Value = "mantener un registro de pedidos  actualizado"

$ws.Range("A11").Value = "US006"
$ws.Range("C11").Value = "Listar los 5 productos mas vendidos "
$ws.Range("D11").Value = "obtener los productos que son mas vendidos "

$ws.Range("A12").Value = "US007"
$ws.Range("B12").Value = "Usuario "
$ws.Range("C12").Value = "Consultar el Stock Actual registrado en el sistema."
$ws.Range("D12").Value = "conocer la situacion general del stock "

$ws.Range("A13").Value = "US008"
$ws.Range("C13").Value = "Actualizar el Stock Actual registrado en el sistema."
$ws.Range("D13").Value = "mantenerse un stock de productos actualizado"

$ws.Range("A14").Value = "US009"
$ws.Range("C14").Value = "Obtener los productos que tienen poco stock minimo "
$ws.Range("D14").Value = "mantener un registro de productos con stock minimo para reponer"

$ws.Range("A15").Value = "US0010"
$ws.Range("C15").Value = "Quiero darme de alta en el sistemas "
$ws.Range("D15").Value = "generar mi cuenta de usuario para ingresar al sistema de gestion"

$ws.Range("A16").Value = "US0011"
$ws.Range("C16").Value = "Quiero resetear mi clave de usuario "
$ws.Range("D16").Value = "ingresar con mi clave nueva al sistema de gestion "

$ws.Range("A17").Value = "US0012"
$ws.Range("C17").Value = "Quiero acceder al sistema de gestion con mi usuario"
$ws.Range("D17").Value = "realizar las gestiones oportunas en el sistema"
# D17 gains an underlined font (reuses the existing underline font already
# present in the workbook) while keeping its left-aligned bordered layout.
$ws.Range("D17").Font.Underline = 2

# ---------------------------------------------------------------------------
# Rows 18-25: brand new backlog entries (previously blank rows).
# ---------------------------------------------------------------------------

$ws.Range("A18").Value = "US0013"
$ws.Range("B18").Value = "Usuario"
$ws.Range("C18").Value = "Como usuario quiero tener un manual de uso de la aplicacion de escritorio."
$ws.Range("D18").Value = " aprender el manejo del aplicativo de manera eficaz."

$ws.Range("A19").Value = "US0014"
$ws.Range("B19").Value = "Usuario"
$ws.Range("C19").Value = "Modificar uno o varios items de una venta  ya registrado en el sistema. "
$ws.Range("D19").Value = "tener un registro correcto de las ventas realizadas."

$ws.Range("A20").Value = "US0015"
$ws.Range("B20").Value = "Usuario"
$ws.Range("C20").Value = "Eliminar una venta que no se ha realizado por algun impedimento "
$ws.Range("D20").Value = "tener un registro correcto de las ventas realizadas. "

$ws.Range("A21").Value = "US0016"
$ws.Range("B21").Value = "Usuario"
$ws.Range("C21").Value = "eliminar stock de un producto determinado"
$ws.Range("D21").Value = "manterse un registro de stock actualizado"

$ws.Range("A22").Value = "US0017"
$ws.Range("B22").Value = "Usuario"
$ws.Range("C22").Value = "modificar datos de producto con stock minimo "
$ws.Range("D22").Value = "manterse un registro de stock actualizado"

$ws.Range("A23").Value = "US0018"
$ws.Range("B23").Value = "Administrador "
$ws.Range("C23").Value = "obtener un listado de ventas realizada "
$ws.Range("D23").Value = "mantener un registro de las ventas realizadas en un periodo de tiempo"

$ws.Range("A24").Value = "US0019"
$ws.Range("B24").Value = "Administrador "
$ws.Range("C24").Value = "obtener un listado de usuarios desactivados del sistema"
$ws.Range("D24").Value = "obtener un listado de usuarios no activos en un periodo de tiempo"

$ws.Range("A25").Value = "US0020"
$ws.Range("B25").Value = "Usuario"
$ws.Range("C25").Value = "quiero poder ejecutar el sistema de gestion de stock en todas las versiones de Windows, desde Windows 95 en adelante"
$ws.Range("D25").Value = "tener el sistema de gestion en funcionamiento en dichas versiones de SO"

# Give the newly populated A18:D25 block the same bordered / filled look
# used by the rest of the table (white fill + thin border on all sides).
$newBlock = $ws.Range("A18:D25")
$newBlock.Borders.LineStyle = 1
$newBlock.Borders.Weight = 2

# C25 uses a distinct small grey Arial note-style font instead of the
# default table font.
$ws.Range("C25").ClearFormats()
$ws.Range("C25").Borders.LineStyle = 1
$ws.Range("C25").Borders.Weight = 2
$ws.Range("C25").Font.Name = "Arial"
$ws.Range("C25").Font.Size = 10
$ws.Range("C25").Font.Color = 2236962

# ---------------------------------------------------------------------------
# Misc formatting / UI state changes.
# ---------------------------------------------------------------------------

# B27 picks up the workbook's underline-font / white-fill note style.
$ws.Range("B27").Font.Underline = 2

# Move the active selection to D25 (matches the author's final cursor spot).
$ws.Range("D25").Select()
